$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 4
$ws.Range("F2").Value = -3
$ws.Range("H2").Value = 46

$ws.Range("E6").Select()
